$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" column header, matching the formatting of the other
# header cells (bold, centered, bordered) by copying an existing header
# cell's format before writing the new text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Populate the "Save" column values for each data row (2-14).
$saveValues = @(0, 0, 0, 1, 1, 0, 0, 0, 0, 0, 1, 0, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
